# Atualização automática de preços de eletricidade
# Updates row 2 of the Spot_PT sheet with the newly scraped day's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 45901
$ws.Range("B2").Value = 35.79
$ws.Range("C2").Value = 18.07
$ws.Range("D2").Value = 22.34
$ws.Range("E2").Value = 17.64
$ws.Range("F2").Value = 15
$ws.Range("G2").Value = 25.21
$ws.Range("H2").Value = 33.23
$ws.Range("I2").Value = 53.18
$ws.Range("J2").Value = 64
$ws.Range("K2").Value = 35.79
$ws.Range("L2").Value = 3.52
$ws.Range("M2").Value = 0.5
$ws.Range("N2").Value = 0.05
$ws.Range("O2").Value = 0
$ws.Range("P2").Value = 0
$ws.Range("Q2").Value = 0
$ws.Range("R2").Value = 0
$ws.Range("S2").Value = 0
$ws.Range("T2").Value = 3.52
$ws.Range("U2").Value = 55
$ws.Range("V2").Value = 80
$ws.Range("W2").Value = 100.01
$ws.Range("X2").Value = 92.8
$ws.Range("Y2").Value = 80
$ws.Range("Z2").Value = 30.65
$ws.Range("AA2").Value = "20h-24h"
$ws.Range("AB2").Value = 88.2
$ws.Range("AC2").Value = "20h-22h"
$ws.Range("AD2").Value = 90
$ws.Range("AE2").Value = "22h-24h"
$ws.Range("AF2").Value = 86.40000000000001
$ws.Range("AG2").Value = "1h-18h"
